# Apply updated cryptocurrency price/volume data to Sheet1.
# Each target cell is forced to Text format before the write (several
# "Price" values look numeric, e.g. 213.60, and Excel would otherwise
# silently coerce them to doubles), then the style is reset back to
# "Normal" so no stray cell-style id is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '27.015.73'
Set-TextValue 'E2' '  -0.55%  '
Set-TextValue 'D3' '1.621.26'
Set-TextValue 'E3' '  -0.96%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '213.60'
Set-TextValue 'E5' '  -1.47%  '
Set-TextValue 'D6' '0.510'
Set-TextValue 'E6' '  -0.81%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'B8' 'Dogecoin'
Set-TextValue 'C8' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue 'D8' '0.0627'
Set-TextValue 'E8' '  +0.39%  '
Set-TextValue 'B9' 'Cardano'
Set-TextValue 'C9' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D9' '0.250'
Set-TextValue 'E9' '  -1.31%  '
Set-TextValue 'D10' '19.92'
Set-TextValue 'E10' '  -0.50%  '
Set-TextValue 'D11' '0.0837'
Set-TextValue 'E11' '  -1.15%  '
Set-TextValue 'D12' '1.847.47'
Set-TextValue 'D13' '1.623.22'
Set-TextValue 'E13' '  -0.76%  '
Set-TextValue 'E14' '  -0.40%  '
Set-TextValue 'E15' '  -0.78%  '
Set-TextValue 'D16' '27.012.67'
Set-TextValue 'E16' '  -0.54%  '
Set-TextValue 'D17' '64.28'
Set-TextValue 'E17' '  -3.11%  '
Set-TextValue 'E18' '  -0.36%  '
Set-TextValue 'D19' '213.17'
Set-TextValue 'E19' '  -1.50%  '
Set-TextValue 'D21' '6.80'
Set-TextValue 'E21' '  -0.30%  '
Set-TextValue 'E22' '  -1.95%  '
Set-TextValue 'D23' '2.35'
Set-TextValue 'E23' '  -7.48%  '
Set-TextValue 'D24' '8.95'
Set-TextValue 'E24' '  -1.82%  '
Set-TextValue 'D25' '147.24'
Set-TextValue 'E25' '  -0.37%  '
Set-TextValue 'D26' '7.48'
Set-TextValue 'E26' '  +1.25%  '
Set-TextValue 'E27' '  +0.11%  '
Set-TextValue 'E28' '  -3.65%  '
Set-TextValue 'D29' '15.48'
Set-TextValue 'E29' '  -1.02%  '
Set-TextValue 'D30' '0.0508'
Set-TextValue 'E30' '  +0.30%  '
Set-TextValue 'E31' '  -1.11%  '
Set-TextValue 'D32' '3.28'
Set-TextValue 'E32' '  -2.73%  '
Set-TextValue 'D33' '0.708'
Set-TextValue 'E33' '  +28.87%  '
Set-TextValue 'E34' '  -0.98%  '
Set-TextValue 'D35' '1.334.41'
Set-TextValue 'E35' '  +2.64%  '
Set-TextValue 'E36' '  -0.70%  '
Set-TextValue 'E37' '  -0.44%  '
Set-TextValue 'E38' '  -0.57%  '
Set-TextValue 'E39' '  -1.77%  '
Set-TextValue 'E40' '  -0.01%  '
Set-TextValue 'D41' '0.796'
Set-TextValue 'E41' '  -1.20%  '
Set-TextValue 'D42' '2.20'
Set-TextValue 'E42' '  -1.74%  '
Set-TextValue 'D43' '5.34'
Set-TextValue 'E43' '  +0.13%  '
Set-TextValue 'D44' '63.81'
Set-TextValue 'E44' '  +3.06%  '
Set-TextValue 'D45' '1.759.31'
Set-TextValue 'E45' '  -0.94%  '
Set-TextValue 'D46' '89.80'
Set-TextValue 'E46' '  -1.64%  '
Set-TextValue 'E47' '  +1.90%  '
Set-TextValue 'E48' '  +21.02%  '
Set-TextValue 'D49' '0.0516'
Set-TextValue 'E49' '  +0.33%  '
Set-TextValue 'D50' '0.0992'
Set-TextValue 'E50' '  +3.82%  '
Set-TextValue 'D51' '7.58'
Set-TextValue 'E51' '  -0.49%  '
